$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.155.53"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.069.91"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.98"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.00"
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.066.93"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.28"
$ws.Range("E10").Value = "  -2.04%  "
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.466"
$ws.Range("E12").Value = "  -3.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000237"
$ws.Range("E13").Value = "  -4.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.60"
$ws.Range("E14").Value = "  -4.09%  "
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.585.42"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.118.60"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.91"
$ws.Range("E18").Value = "  -3.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.073.19"
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.40"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "483.37"
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.683"
$ws.Range("E22").Value = "  -3.66%  "
$ws.Range("E23").Value = "  -3.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.28"
$ws.Range("E24").Value = "  -1.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.59"
$ws.Range("E25").Value = "  -4.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.19"
$ws.Range("E26").Value = "  -3.44%  "
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.85"
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("E30").Value = "  -5.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.62"
$ws.Range("E32").Value = "  -2.95%  "
$ws.Range("E33").Value = "  -4.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0898"
$ws.Range("E34").Value = "  -5.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "47.27"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("E37").Value = "  -3.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.54"
$ws.Range("E38").Value = "  -4.95%  "
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.94"
$ws.Range("E40").Value = "  -5.65%  "
$ws.Range("E41").Value = "  -4.21%  "
$ws.Range("E42").Value = "  -5.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.762.99"
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "134.53"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0341"
$ws.Range("E46").Value = "  -3.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "360.93"
$ws.Range("E47").Value = "  -5.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.14"
$ws.Range("E49").Value = "  -2.63%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.106"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.13"
$ws.Range("E51").Value = "  -2.84%  "
